# Add two new machine records (Machine 30 and Machine 31) to the master
# machine data table, continuing from the last existing row (row 30).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 31 - Machine 30
$ws.Range("A31").Value = 10030
$ws.Range("B31").Value = "Machine 30"
$ws.Range("C31").Value = "70-5A-0F-8C-01-39"
$ws.Range("D31").Value = "FB5962911663"
$ws.Range("E31").Value = "192.168.0.356"
$ws.Range("F31").Value = 1001
$ws.Range("G31").Value = "eng"
$ws.Range("H31").Value = $true
$ws.Range("I31").Value = "superadmin"
$ws.Range("J31").Value = "now()"
$ws.Range("K31").Value = "now()"

# Row 32 - Machine 31
# NOTE: cell values are written in a specific order so that new shared
# strings are interned in the same order as the target workbook (the
# ip_address string is first used before the mac_address string).
$ws.Range("A32").Value = 10031
$ws.Range("B32").Value = "Machine 31"
$ws.Range("E32").Value = "192.168.0.357"
$ws.Range("C32").Value = "58-20-B1-DA-F3-FB"
$ws.Range("D32").Value = "FB5962911663"
$ws.Range("F32").Value = 1001
$ws.Range("G32").Value = "eng"
$ws.Range("H32").Value = $true
$ws.Range("I32").Value = "superadmin"
$ws.Range("J32").Value = "now()"
$ws.Range("K32").Value = "now()"

# Move the "new row" selection marker down to reflect the fact that data
# now ends at row 32 (selection of the empty area below the table).
$ws.Range("L1:XFD1048576").Select()
